$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.627
$ws.Range("D7").Value = -7.483000000000001
$ws.Range("A8").Value = -22.188
$ws.Range("A10").Value = -21.846
$ws.Range("A12").Value = -21.589
$ws.Range("D15").Value = -8.147000000000002
$ws.Range("A18").Value = -22.166
$ws.Range("D18").Value = -8.4
$ws.Range("E18").Value = 16.366
$ws.Range("E19").Value = 16.538
$ws.Range("D20").Value = -7.57
$ws.Range("E27").Value = 16.38
$ws.Range("D29").Value = -7.292
$ws.Range("D30").Value = -7.048999999999999
$ws.Range("D31").Value = -7.841999999999999
$ws.Range("E31").Value = 16.956
$ws.Range("A37").Value = -20.21700000000001
$ws.Range("E38").Value = 16.726
$ws.Range("D40").Value = -7.780999999999999
$ws.Range("E42").Value = 16.665
$ws.Range("E44").Value = 16.788
$ws.Range("E47").Value = 16.471
$ws.Range("D50").Value = -8.105
$ws.Range("A55").Value = -22.283
$ws.Range("E58").Value = 16.448
$ws.Range("E65").Value = 17.301
$ws.Range("A68").Value = -21.534
$ws.Range("D68").Value = -6.778
$ws.Range("E73").Value = 16.635
$ws.Range("D76").Value = -7.753000000000002
$ws.Range("A77").Value = -20.963
$ws.Range("A78").Value = -20.281
$ws.Range("A81").Value = -21.82
$ws.Range("A82").Value = -22.261
$ws.Range("D87").Value = -8.261999999999999
$ws.Range("D88").Value = -8.293000000000001
$ws.Range("E90").Value = 16.492
$ws.Range("E94").Value = 17.686
$ws.Range("E95").Value = 17.397
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.242000000000001
$ws.Range("D101").Value = -7.616
$ws.Range("E101").Value = 16.44
$ws.Range("D102").Value = -8.142999999999999
